$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (15,1,... and 20,1,...) to append starting at row 14, columns A-Q
$data = @(
    @(15,1,26.886,402.026,2.14549304,18234644.416,0.11276818,0.20143332,5.54187442,4.389023930457005,65.41304960110202,0.338025257771874,6521921.218650043,0.02163876988212646,0.04825498426058783,2.074995770643459,0.15),
    @(15,1,43.792,642.004,1.33977116,27555683.248,0.10412572,0.30250084,13.51638904,7.028765402698567,99.66933297595476,0.1979011030465178,9487486.102261616,0.01570389910103701,0.05884070019376882,4.544064589209102,0.85),
    @(20,1,19.88,395.364,2.1927889,31197066.008,0.1453059,0.14541756,3.01728986,3.684664639694368,72.45569100020117,0.3715919961821612,11093264.63146623,0.02876172036068448,0.04663943911979666,1.68096522936903,0.15),
    @(20,1,33.412,647.138,1.33553914,48918644.056,0.133752,0.22370094,7.72350346,5.919149318720945,109.8992203487656,0.2174971694705837,15225068.91508796,0.01751580343403864,0.05141625244363006,3.214486784118169,0.85)
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
